$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 61, shifting old row 61 (and below) down.
$ws.Rows.Item(61).Insert()

# Fill the new row 61 with the new weekly price data.
$ws.Cells.Item(61, 1).Value = 1
$ws.Cells.Item(61, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(61, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(61, 4).Value = 44568
$ws.Cells.Item(61, 5).Value = 15
$ws.Cells.Item(61, 6).Value = 100112038
$ws.Cells.Item(61, 7).Value = "Cebollín baby"
$ws.Cells.Item(61, 8).Value = "Sin especificar"
$ws.Cells.Item(61, 9).Value = "Primera"
$ws.Cells.Item(61, 10).Value = 300
$ws.Cells.Item(61, 11).Value = 5000
$ws.Cells.Item(61, 12).Value = 5500
$ws.Cells.Item(61, 13).Value = 5250
$ws.Cells.Item(61, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(61, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(61, 16).Value = 2625
$ws.Cells.Item(61, 17).Value = 2
$ws.Cells.Item(61, 18).Value = "Hortaliza"
